$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.250631
$ws.Range("H2").Value = 0.751893
$ws.Range("I2").Value = 0.2648339568266264
$ws.Range("J2").Value = 0.2648339568266264
$ws.Range("M2").Value = 36.923013
$ws.Range("N2").Value = 110.769039
$ws.Range("O2").Value = 0.7437819354528793
$ws.Range("P2").Value = 0.7437819354528794
$ws.Range("Q2").Value = 9.254051671202999
$ws.Range("R2").Value = 83.286465040827
$ws.Range("S2").Value = 0.1969787129821524
$ws.Range("T2").Value = 0.1969787129821525

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.250631
$ws.Range("H3").Value = 0.751893
$ws.Range("I3").Value = 0.2648339568266264
$ws.Range("J3").Value = 0.2648339568266264
$ws.Range("O3").Value = 0.1364233939221953
$ws.Range("P3").Value = 0.1364233939221953
$ws.Range("Q3").Value = 1.697364612315
$ws.Range("R3").Value = 15.276281510835
$ws.Range("S3").Value = 0.03612954721613252
$ws.Range("T3").Value = 0.03612954721613252

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.250631
$ws.Range("H4").Value = 0.751893
$ws.Range("I4").Value = 0.2648339568266264
$ws.Range("J4").Value = 0.2648339568266264
$ws.Range("O4").Value = 0.1197946706249253
$ws.Range("P4").Value = 0.1197946706249254
$ws.Range("Q4").Value = 1.490471896474333
$ws.Range("R4").Value = 13.414247068269
$ws.Range("S4").Value = 0.0317256966283414
$ws.Range("T4").Value = 0.0317256966283414

$ws.Range("I5").Value = 0.2480790641859371
$ws.Range("J5").Value = 0.2480790641859371
$ws.Range("M5").Value = 36.923013
$ws.Range("N5").Value = 110.769039
$ws.Range("O5").Value = 0.7437819354528793
$ws.Range("P5").Value = 0.7437819354528794
$ws.Range("Q5").Value = 8.668588069404001
$ws.Range("R5").Value = 78.017292624636
$ws.Range("S5").Value = 0.1845167265055554
$ws.Range("T5").Value = 0.1845167265055554

$ws.Range("I6").Value = 0.2480790641859371
$ws.Range("J6").Value = 0.2480790641859371
$ws.Range("O6").Value = 0.1364233939221953
$ws.Range("P6").Value = 0.1364233939221953
$ws.Range("S6").Value = 0.03384378789728767
$ws.Range("T6").Value = 0.03384378789728767

$ws.Range("I7").Value = 0.2480790641859371
$ws.Range("J7").Value = 0.2480790641859371
$ws.Range("O7").Value = 0.1197946706249253
$ws.Range("P7").Value = 0.1197946706249254
$ws.Range("S7").Value = 0.02971854978309405
$ws.Range("T7").Value = 0.02971854978309405

$ws.Range("G8").Value = 0.4609646666666666
$ws.Range("I8").Value = 0.4870869789874365
$ws.Range("J8").Value = 0.4870869789874365
$ws.Range("M8").Value = 36.923013
$ws.Range("N8").Value = 110.769039
$ws.Range("O8").Value = 0.7437819354528793
$ws.Range("P8").Value = 0.7437819354528794
$ws.Range("Q8").Value = 17.020204379874
$ws.Range("R8").Value = 153.181839418866
$ws.Range("S8").Value = 0.3622864959651715
$ws.Range("T8").Value = 0.3622864959651715

$ws.Range("G9").Value = 0.4609646666666666
$ws.Range("I9").Value = 0.4870869789874365
$ws.Range("J9").Value = 0.4870869789874365
$ws.Range("O9").Value = 0.1364233939221953
$ws.Range("P9").Value = 0.1364233939221953
$ws.Range("R9").Value = 28.09638877292999
$ws.Range("S9").Value = 0.06645005880877514
$ws.Range("T9").Value = 0.06645005880877514

$ws.Range("G10").Value = 0.4609646666666666
$ws.Range("I10").Value = 0.4870869789874365
$ws.Range("J10").Value = 0.4870869789874365
$ws.Range("O10").Value = 0.1197946706249253
$ws.Range("P10").Value = 0.1197946706249254
$ws.Range("S10").Value = 0.05835042421348989
$ws.Range("T10").Value = 0.0583504242134899

Write-Output "Applied updates to Sema3a-Plxna2 sheet"